$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 13331.8
$ws.Range("I7").Value = 2221.3333
$ws.Range("J7").Value = 29997.5
$ws.Range("K7").Value = 2221.3333
$ws.Range("L7").Value = 29997.5
$ws.Range("M7").Value = -2109.3333

$ws.Range("H14").Value = 13331.8
$ws.Range("I14").Value = 2221.3333
$ws.Range("J14").Value = 29997.5
$ws.Range("K14").Value = 2221.3333
$ws.Range("L14").Value = 29997.5
$ws.Range("M14").Value = -2030.3333

$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = $null

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = $null

$ws.Range("H123").Value = 180000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 180000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 180000
$ws.Range("N123").Value = -189800

$ws.Range("H125").Value = 2395
$ws.Range("I125").Value = 2395
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 21555
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -19095
$ws.Range("N125").Value = $null

$ws.Range("H131").Value = 7996.3335
$ws.Range("I131").Value = 7994.5
$ws.Range("J131").Value = 8000
$ws.Range("K131").Value = 23983.5
$ws.Range("L131").Value = 24000
$ws.Range("M131").Value = -18943.5
$ws.Range("N131").Value = -34080

$ws.Range("H138").Value = 2535.3276
$ws.Range("I138").Value = 1759.0625
$ws.Range("J138").Value = 2831.0476
$ws.Range("K138").Value = 5277.1875
$ws.Range("L138").Value = 8493.1428
$ws.Range("M138").Value = -137.1875
$ws.Range("N138").Value = -18773.1428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 180
$ws.Range("I4").Value = 180
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 180
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -64
$ws.Range("N4").Value = $null

$ws.Range("H32").Value = 4886.55
$ws.Range("I32").Value = 4886.55
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4886.55
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -4599.55

$ws.Range("H41").Value = 2056
$ws.Range("I41").Value = 2056
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 2056
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -1642

$ws.Range("H61").Value = 1824.4166
$ws.Range("I61").Value = 1299.875
$ws.Range("J61").Value = 2873.5
$ws.Range("K61").Value = 1299.875
$ws.Range("L61").Value = 2873.5
$ws.Range("M61").Value = -1087.875
$ws.Range("N61").Value = -3297.5

$ws.Range("H136").Value = 1824.4166
$ws.Range("I136").Value = 1299.875
$ws.Range("J136").Value = 2873.5
$ws.Range("K136").Value = 3899.625
$ws.Range("L136").Value = 8620.5
$ws.Range("M136").Value = -1349.625
$ws.Range("N136").Value = -13720.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2000
$ws.Range("I99").Value = 2000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -502
$ws.Range("N99").Value = $null

$ws.Range("H105").Value = 5665
$ws.Range("I105").Value = 5665
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 5665
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -3918

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").Value = $null

$ws.Range("H58").Value = 2462.3333
$ws.Range("I58").Value = 1924.5
$ws.Range("J58").Value = 2731.25
$ws.Range("K58").Value = 1924.5
$ws.Range("L58").Value = 2731.25
$ws.Range("M58").Value = -1721.5

$ws.Range("H74").Value = 49976
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 49976
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 49976
$ws.Range("N74").Value = -51724

$ws.Range("H77").Value = 49976
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 49976
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 149928
$ws.Range("N77").Value = -158664

$ws.Range("H136").Value = 2462.3333
$ws.Range("I136").Value = 1924.5
$ws.Range("J136").Value = 2731.25
$ws.Range("K136").Value = 5773.5
$ws.Range("L136").Value = 8193.75
$ws.Range("M136").Value = -3223.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 220908240
$ws.Range("I4").Value = 131098480
$ws.Range("J4").Value = 400527780
$ws.Range("K4").Value = 393295440
$ws.Range("L4").Value = 1201583340
$ws.Range("M4").Value = -393295328

$ws.Range("H12").Value = 213.2
$ws.Range("I12").Value = 231
$ws.Range("J12").Value = 201.33333
$ws.Range("K12").Value = 693
$ws.Range("L12").Value = 603.99999
$ws.Range("M12").Value = -520
$ws.Range("N12").Value = -949.99999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 65.09999999999999
$ws.Range("I2").Value = 67.625
$ws.Range("J2").Value = 55
$ws.Range("K2").Value = 67.625
$ws.Range("L2").Value = 55
$ws.Range("M2").Value = 45.375

$ws.Range("H70").Value = 6000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 6000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 6000
$ws.Range("N70").Value = -6540

$ws.Range("H73").Value = 6000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 6000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 6000
$ws.Range("N73").Value = -7872

$ws.Range("H97").Value = 5000
$ws.Range("I97").Value = 5000
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 5000
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -4504

$ws.Range("H132").Value = 2129.35
$ws.Range("I132").Value = 1080.1111
$ws.Range("J132").Value = 2987.818
$ws.Range("K132").Value = 3240.3333
$ws.Range("L132").Value = 8963.454000000002
$ws.Range("M132").Value = -710.3333000000002
$ws.Range("N132").Value = -14023.454

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = $null

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = $null

$ws.Range("H43").Value = 189999.2
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 189999.2
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 189999.2
$ws.Range("N43").Value = -190385.2

$ws.Range("H61").Value = 698.3333
$ws.Range("I61").Value = 650
$ws.Range("J61").Value = 795
$ws.Range("K61").Value = 650
$ws.Range("L61").Value = 795
$ws.Range("M61").Value = -448

$ws.Range("H113").Value = 698.3333
$ws.Range("I113").Value = 650
$ws.Range("J113").Value = 795
$ws.Range("K113").Value = 650
$ws.Range("L113").Value = 795
$ws.Range("M113").Value = 1520

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 30656.5
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 30656.5
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 30656.5
$ws.Range("N45").Value = -31638.5

$ws.Range("H100").Value = 2300
$ws.Range("I100").Value = 2300
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 4600
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -4059
$ws.Range("N100").Value = $null
